$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix minor errors on map: correct the "Farmstand Local Foods" entries to note
# that it is a distributor, and fix the Vashon-Maury Island Land Trust row to
# map to the correct farm name (Matsuda Farm).
$ws.Range("B25").Value = "Farmstand Local Foods (distributor)"
$ws.Range("B26").Value = "Farmstand Local Foods (distributor)"
$ws.Range("B98").Value = "Matsuda Farm"

# Update the view to reflect where the author was working when they saved.
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("B27").Select()
